$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.23492266666667
$ws.Range("H2").Value = 45.704768
$ws.Range("I2").Value = 0.245239930755761
$ws.Range("J2").Value = 0.245239930755761
$ws.Range("M2").Value = 0.9601406666666668
$ws.Range("N2").Value = 2.880422
$ws.Range("O2").Value = 0.004534252661098308
$ws.Range("P2").Value = 0.004534252661098308
$ws.Range("Q2").Value = 14.62766880578845
$ws.Range("R2").Value = 131.649019252096
$ws.Range("S2").Value = 0.001111979808636874
$ws.Range("T2").Value = 0.001111979808636874
$ws.Range("G3").Value = 15.23492266666667
$ws.Range("H3").Value = 45.704768
$ws.Range("I3").Value = 0.245239930755761
$ws.Range("J3").Value = 0.245239930755761
$ws.Range("O3").Value = 0.01855104847259158
$ws.Range("P3").Value = 0.01855104847259158
$ws.Range("Q3").Value = 59.84637675472356
$ws.Range("R3").Value = 538.617390792512
$ws.Range("S3").Value = 0.004549457842865125
$ws.Range("T3").Value = 0.004549457842865125
$ws.Range("G4").Value = 15.23492266666667
$ws.Range("H4").Value = 45.704768
$ws.Range("I4").Value = 0.245239930755761
$ws.Range("J4").Value = 0.245239930755761
$ws.Range("M4").Value = 1.316544333333334
$ws.Range("N4").Value = 3.949633
$ws.Range("O4").Value = 0.006217364657196653
$ws.Range("P4").Value = 0.006217364657196652
$ws.Range("Q4").Value = 20.05745110557156
$ws.Range("R4").Value = 180.517059950144
$ws.Range("S4").Value = 0.001524746078014223
$ws.Range("T4").Value = 0.001524746078014223
$ws.Range("G5").Value = 15.23492266666667
$ws.Range("H5").Value = 45.704768
$ws.Range("I5").Value = 0.245239930755761
$ws.Range("J5").Value = 0.245239930755761
$ws.Range("M5").Value = 205.547872
$ws.Range("N5").Value = 616.643616
$ws.Range("O5").Value = 0.9706973342091134
$ws.Range("P5").Value = 0.9706973342091134
$ws.Range("Q5").Value = 3131.505934217899
$ws.Range("R5").Value = 28183.55340796109
$ws.Range("S5").Value = 0.2380537470262448
$ws.Range("T5").Value = 0.2380537470262448
$ws.Range("I6").Value = 0.623546923900845
$ws.Range("J6").Value = 0.6235469239008449
$ws.Range("M6").Value = 0.9601406666666668
$ws.Range("N6").Value = 2.880422
$ws.Range("O6").Value = 0.004534252661098308
$ws.Range("P6").Value = 0.004534252661098308
$ws.Range("Q6").Value = 37.19230330713778
$ws.Range("R6").Value = 334.73072976424
$ws.Range("S6").Value = 0.002827319299017071
$ws.Range("T6").Value = 0.00282731929901707
$ws.Range("I7").Value = 0.623546923900845
$ws.Range("J7").Value = 0.6235469239008449
$ws.Range("O7").Value = 0.01855104847259158
$ws.Range("P7").Value = 0.01855104847259158
$ws.Range("S7").Value = 0.01156744921021995
$ws.Range("T7").Value = 0.01156744921021994
$ws.Range("I8").Value = 0.623546923900845
$ws.Range("J8").Value = 0.6235469239008449
$ws.Range("M8").Value = 1.316544333333334
$ws.Range("N8").Value = 3.949633
$ws.Range("O8").Value = 0.006217364657196653
$ws.Range("P8").Value = 0.006217364657196652
$ws.Range("Q8").Value = 50.99806503626223
$ws.Range("R8").Value = 458.9825853263601
$ws.Range("S8").Value = 0.003876818606764804
$ws.Range("T8").Value = 0.003876818606764803
$ws.Range("I9").Value = 0.623546923900845
$ws.Range("J9").Value = 0.6235469239008449
$ws.Range("M9").Value = 205.547872
$ws.Range("N9").Value = 616.643616
$ws.Range("O9").Value = 0.9706973342091134
$ws.Range("P9").Value = 0.9706973342091134
$ws.Range("Q9").Value = 7962.165404472747
$ws.Range("R9").Value = 71659.48864025471
$ws.Range("S9").Value = 0.6052753367848431
$ws.Range("T9").Value = 0.605275336784843
$ws.Range("G10").Value = 6.113134
$ws.Range("H10").Value = 18.339402
$ws.Range("I10").Value = 0.09840447448682081
$ws.Range("J10").Value = 0.09840447448682083
$ws.Range("M10").Value = 0.9601406666666668
$ws.Range("N10").Value = 2.880422
$ws.Range("O10").Value = 0.004534252661098308
$ws.Range("P10").Value = 0.004534252661098308
$ws.Range("Q10").Value = 5.869468554182667
$ws.Range("R10").Value = 52.825216987644
$ws.Range("S10").Value = 0.0004461907503058478
$ws.Range("T10").Value = 0.0004461907503058479
$ws.Range("G11").Value = 6.113134
$ws.Range("H11").Value = 18.339402
$ws.Range("I11").Value = 0.09840447448682081
$ws.Range("J11").Value = 0.09840447448682083
$ws.Range("O11").Value = 0.01855104847259158
$ws.Range("P11").Value = 0.01855104847259158
$ws.Range("Q11").Value = 24.01383508933533
$ws.Range("R11").Value = 216.124515804018
$ws.Range("S11").Value = 0.001825506176124914
$ws.Range("T11").Value = 0.001825506176124914
$ws.Range("G12").Value = 6.113134
$ws.Range("H12").Value = 18.339402
$ws.Range("I12").Value = 0.09840447448682081
$ws.Range("J12").Value = 0.09840447448682083
$ws.Range("M12").Value = 1.316544333333334
$ws.Range("N12").Value = 3.949633
$ws.Range("O12").Value = 0.006217364657196653
$ws.Range("P12").Value = 0.006217364657196652
$ws.Range("Q12").Value = 8.048211926607333
$ws.Range("R12").Value = 72.433907339466
$ws.Range("S12").Value = 0.0006118165017843694
$ws.Range("T12").Value = 0.0006118165017843694
$ws.Range("G13").Value = 6.113134
$ws.Range("H13").Value = 18.339402
$ws.Range("I13").Value = 0.09840447448682081
$ws.Range("J13").Value = 0.09840447448682083
$ws.Range("M13").Value = 205.547872
$ws.Range("N13").Value = 616.643616
$ws.Range("O13").Value = 0.9706973342091134
$ws.Range("P13").Value = 0.9706973342091134
$ws.Range("Q13").Value = 1256.541684950848
$ws.Range("R13").Value = 11308.87516455763
$ws.Range("S13").Value = 0.09552096105860568
$ws.Range("T13").Value = 0.09552096105860569
$ws.Range("G14").Value = 2.038157333333333
$ws.Range("H14").Value = 6.114472
$ws.Range("I14").Value = 0.0328086708565732
$ws.Range("J14").Value = 0.0328086708565732
$ws.Range("M14").Value = 0.9601406666666668
$ws.Range("N14").Value = 2.880422
$ws.Range("O14").Value = 0.004534252661098308
$ws.Range("P14").Value = 0.004534252661098308
$ws.Range("Q14").Value = 1.956917740798222
$ws.Range("R14").Value = 17.612259667184
$ws.Range("S14").Value = 0.0001487628031385155
$ws.Range("T14").Value = 0.0001487628031385155
$ws.Range("G15").Value = 2.038157333333333
$ws.Range("H15").Value = 6.114472
$ws.Range("I15").Value = 0.0328086708565732
$ws.Range("J15").Value = 0.0328086708565732
$ws.Range("O15").Value = 0.01855104847259158
$ws.Range("P15").Value = 0.01855104847259158
$ws.Range("Q15").Value = 8.006363689849778
$ws.Range("R15").Value = 72.057273208648
$ws.Range("S15").Value = 0.0006086352433815921
$ws.Range("T15").Value = 0.0006086352433815921
$ws.Range("G16").Value = 2.038157333333333
$ws.Range("H16").Value = 6.114472
$ws.Range("I16").Value = 0.0328086708565732
$ws.Range("J16").Value = 0.0328086708565732
$ws.Range("M16").Value = 1.316544333333334
$ws.Range("N16").Value = 3.949633
$ws.Range("O16").Value = 0.006217364657196653
$ws.Range("P16").Value = 0.006217364657196652
$ws.Range("Q16").Value = 2.683324487641778
$ws.Range("R16").Value = 24.149920388776
$ws.Range("S16").Value = 0.000203983470633256
$ws.Range("T16").Value = 0.000203983470633256
$ws.Range("G17").Value = 2.038157333333333
$ws.Range("H17").Value = 6.114472
$ws.Range("I17").Value = 0.0328086708565732
$ws.Range("J17").Value = 0.0328086708565732
$ws.Range("M17").Value = 205.547872
$ws.Range("N17").Value = 616.643616
$ws.Range("O17").Value = 0.9706973342091134
$ws.Range("P17").Value = 0.9706973342091134
$ws.Range("Q17").Value = 418.9389026678613
$ws.Range("R17").Value = 3770.450124010752
$ws.Range("S17").Value = 0.03184728933941983
$ws.Range("T17").Value = 0.03184728933941983
